# Update Name of Algo
# Apply updated imputed values in the KNN result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.715999999999999
$ws.Range("A12").Value = -21.626
$ws.Range("D14").Value = -7.716000000000001
$ws.Range("D26").Value = -8.000999999999999
$ws.Range("D31").Value = -8.219000000000001
$ws.Range("A32").Value = -21.648
$ws.Range("D35").Value = -7.939
$ws.Range("A36").Value = -20.339
$ws.Range("D37").Value = -7.741
$ws.Range("A38").Value = -19.741
$ws.Range("D45").Value = -7.539
$ws.Range("A46").Value = -21.829
$ws.Range("A54").Value = -21.869
$ws.Range("A55").Value = -22.135
$ws.Range("D57").Value = -8.271000000000001
$ws.Range("A67").Value = -21.565
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.445
$ws.Range("A91").Value = -21.522
$ws.Range("A99").Value = -20.43
$ws.Range("D100").Value = -8.280000000000001
$ws.Range("D102").Value = -7.752
